# Weekly fruit/hortaliza update: a new price observation for Choclo
# (Vega Monumental Concepción) is inserted ahead of the existing row 157,
# pushing the previous rows 157-168 down to 158-169.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 157; existing rows 157:168 shift to 158:169.
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(157, 1).Value = 11
$ws.Cells.Item(157, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(157, 3).Value = "Bíobío"
$ws.Cells.Item(157, 4).Value = 45013
$ws.Cells.Item(157, 5).Value = 8
$ws.Cells.Item(157, 6).Value = 100112024
$ws.Cells.Item(157, 7).Value = "Choclo"
$ws.Cells.Item(157, 8).Value = "Choclero"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 16000
$ws.Cells.Item(157, 11).Value = 350
$ws.Cells.Item(157, 12).Value = 400
$ws.Cells.Item(157, 13).Value = 375
$ws.Cells.Item(157, 14).Value = "$/unidad"
$ws.Cells.Item(157, 15).Value = "Región Metropolitana"
$ws.Cells.Item(157, 16).Value = 375
$ws.Cells.Item(157, 17).Value = 1
$ws.Cells.Item(157, 18).Value = "Hortaliza"
